# "Generate Report for Handoff"
# The localization status report is refreshed: the file
# 81015c87-6f0f-49eb-bec6-73ae877e7b2e.md has just been handed off again,
# so its "Latest Handoff Datetime" (zh-cn / de-de sheets) and
# "Latest HO Xliff Generate Date" (Overview sheet) timestamps are updated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 5 (81015c87-...md), column G "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-10-26 07:14:09"

# --- zh-cn sheet: row 5 (81015c87-...md), column H "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-10-26 07:13:57"

# --- de-de sheet: row 5 (81015c87-...md), column H "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-10-26 07:14:09"
